# "Generate Report for Handoff"
#
# The localization status report is regenerated: the "7f919e97-...md" and
# "8022417e-...md" files (rows 4-5) finish translation and the
# "2339ca68-...md" / "9cdb3ef7-...md" files (rows 6-7) get a fresh handoff,
# so for every language sheet (zh-cn, de-de) rows 4-7:
#   - Priority flips from "low" to "ht" (handoff/translate priority)
#   - Latest Handoff Datetime (column H) is bumped to the new generation time
# The Overview sheet's "Latest HO Xliff Generate Date" column (G) mirrors
# the same timestamp for those rows.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newHandoffTimeZhCn = "2016-10-19 12:28:01"
$newHandoffTimeDeDe = "2016-10-19 12:28:13"

for ($row = 4; $row -le 7; $row++) {
    # Priority column (E) for both language sheets
    $wsZhCn.Cells.Item($row, 5).Value = "ht"
    $wsDeDe.Cells.Item($row, 5).Value = "ht"

    # Latest Handoff Datetime column (H) for both language sheets
    $wsZhCn.Cells.Item($row, 8).Value = $newHandoffTimeZhCn
    $wsDeDe.Cells.Item($row, 8).Value = $newHandoffTimeDeDe

    # Overview sheet's Latest HO Xliff Generate Date (G) mirrors de-de's
    # handoff timestamp for these rows (same underlying shared value).
    $wsOverview.Cells.Item($row, 7).Value = $newHandoffTimeDeDe
}
